# Updates cryptos list values to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.656.14"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.222.99"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.42"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.76"
$ws.Range("E7").Value = "  -7.62%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.401"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.80"
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0892"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").Value = "2.554.91"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.38"
$ws.Range("E14").Value = "  -5.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.38"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.59"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.796"
$ws.Range("E17").Value = "  -4.87%  "
$ws.Range("D18").Value = "2.239.49"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "41.585.27"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.20"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.09"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.32"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.35"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.40"
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.140"
$ws.Range("E29").Value = "  -3.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.80"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.40"
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.54"
$ws.Range("E32").Value = "  -9.68%  "
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.96"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.65"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0650"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.51"
$ws.Range("E37").Value = "  -8.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.38"
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.57"
$ws.Range("E39").Value = "  -8.57%  "
$ws.Range("B40").Value = "BinanceUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("B41").Value = "TerraClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.000235"
$ws.Range("E41").Value = "  -5.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0237"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.60"
$ws.Range("E43").Value = "  -3.42%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0964"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.28"
$ws.Range("E45").Value = "  -4.49%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.21"
$ws.Range("E46").Value = "  -2.76%  "
$ws.Range("E47").Value = "  -9.37%  "
$ws.Range("D48").Value = "1.472.13"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.47"
$ws.Range("E49").Value = "  -7.93%  "
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.07"
$ws.Range("E51").Value = "  -4.29%  "
